# Apply updated evaluation_time values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.01814532279968262
$ws.Range("C2").Value = 0.03399190902709961
$ws.Range("D2").Value = 0.008392953872680664
$ws.Range("E2").Value = 0.01980352401733398
$ws.Range("F2").Value = 0.003587245941162109
$ws.Range("G2").Value = 0.06248388290405273
$ws.Range("H2").Value = 0.01741399765014649
$ws.Range("I2").Value = 0.02348895072937012
$ws.Range("J2").Value = 0.01260161399841309
$ws.Range("K2").Value = 0.02331080436706543
$ws.Range("L2").Value = 0.005994272232055664
$ws.Range("M2").Value = 0.02220544815063476
$ws.Range("B3").Value = 0.06900496482849121
$ws.Range("C3").Value = 0.02339510917663574
$ws.Range("D3").Value = 0.01394104957580566
$ws.Range("E3").Value = 0.009999418258666992
$ws.Range("F3").Value = 0.006599760055541993
$ws.Range("G3").Value = 0.007600641250610352
$ws.Range("H3").Value = 0.1036828994750977
$ws.Range("I3").Value = 0.03320541381835938
$ws.Range("J3").Value = 0.08026003837585449
$ws.Range("K3").Value = 0.02751049995422363
$ws.Range("L3").Value = 0.02319622039794922
$ws.Range("M3").Value = 0.01260318756103516
$ws.Range("B4").Value = 0.02921066284179687
$ws.Range("C4").Value = 0.01479673385620117
$ws.Range("D4").Value = 0.01159276962280273
$ws.Range("E4").Value = 0.009012937545776367
$ws.Range("F4").Value = 0.06995935440063476
$ws.Range("G4").Value = 0.008006906509399414
$ws.Range("H4").Value = 0.01999883651733398
$ws.Range("I4").Value = 0.01385564804077149
$ws.Range("J4").Value = 0.01618986129760742
$ws.Range("K4").Value = 0.0126190185546875
$ws.Range("L4").Value = 0.03120908737182617
$ws.Range("M4").Value = 0.008501768112182617
$ws.Range("B5").Value = 0.01721506118774414
$ws.Range("C5").Value = 0.01379256248474121
$ws.Range("D5").Value = 0.01432771682739258
$ws.Range("E5").Value = 0.01284389495849609
$ws.Range("H5").Value = 0.0125917911529541
$ws.Range("I5").Value = 0.01320748329162598
$ws.Range("J5").Value = 0.009984683990478516
$ws.Range("K5").Value = 0.0116426944732666
$ws.Range("B6").Value = 0.3273634910583496
$ws.Range("C6").Value = 0.05518450736999512
$ws.Range("D6").Value = 0.3017604827880859
$ws.Range("E6").Value = 0.05480003356933594
$ws.Range("F6").Value = 0.1248418807983398
$ws.Range("G6").Value = 0.01859469413757324
$ws.Range("H6").Value = 0.4201132297515869
$ws.Range("I6").Value = 0.06013646125793457
$ws.Range("J6").Value = 0.2652891159057617
$ws.Range("K6").Value = 0.04642405509948731
$ws.Range("L6").Value = 0.1194005489349365
$ws.Range("M6").Value = 0.01959366798400879
